# Update the repayment strategy value on the "ProductLoanInput" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# B17 ("repaymentstrategy" row) changes from "RBI (India)" to the new scenario value.
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Move/refresh the active selection to the edited cell, as recorded in the sheet view.
$ws.Activate()
$ws.Range("B17").Select()
